$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift existing rows 93..113 down to 94..114 (bottom-up) to make room
# for the newly inserted row 93 (date 2021-02-08).
for ($r = 113; $r -ge 93; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value()
}

# Step 2: give the two brand-new rows (113 used to be the last row, now we add
# rows 114 and 115) the same date-column formatting as the rest of column A,
# by copying formats from A112 (an existing, correctly formatted date cell).
$ws.Range("A112").Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: write the new/updated values for rows 90-115 (columns A-D), matching
# the refreshed 7-day rolling sums and the two additional days of data.
$ws.Cells.Item(90, 1).Value = 44232
$ws.Cells.Item(90, 2).Value = 3
$ws.Cells.Item(90, 3).Value = 13
$ws.Cells.Item(90, 4).Value = 406.3769928102532
$ws.Cells.Item(91, 1).Value = 44233
$ws.Cells.Item(91, 2).Value = 3
$ws.Cells.Item(91, 3).Value = 14
$ws.Cells.Item(91, 4).Value = 437.636761487965
$ws.Cells.Item(92, 1).Value = 44234
$ws.Cells.Item(92, 2).Value = 2
$ws.Cells.Item(92, 3).Value = 14
$ws.Cells.Item(92, 4).Value = 437.636761487965
$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 3
$ws.Cells.Item(93, 3).Value = 13
$ws.Cells.Item(93, 4).Value = 406.3769928102532
$ws.Cells.Item(94, 1).Value = 44236
$ws.Cells.Item(94, 2).Value = 2
$ws.Cells.Item(94, 3).Value = 12
$ws.Cells.Item(94, 4).Value = 375.1172241325414
$ws.Cells.Item(95, 1).Value = 44237
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(95, 3).Value = 11
$ws.Cells.Item(95, 4).Value = 343.8574554548296
$ws.Cells.Item(96, 1).Value = 44238
$ws.Cells.Item(96, 2).Value = 0
$ws.Cells.Item(96, 3).Value = 9
$ws.Cells.Item(96, 4).Value = 281.3379180994061
$ws.Cells.Item(97, 1).Value = 44239
$ws.Cells.Item(97, 2).Value = 2
$ws.Cells.Item(97, 3).Value = 6
$ws.Cells.Item(97, 4).Value = 187.5586120662707
$ws.Cells.Item(98, 1).Value = 44240
$ws.Cells.Item(98, 2).Value = 2
$ws.Cells.Item(98, 3).Value = 4
$ws.Cells.Item(98, 4).Value = 125.0390747108471
$ws.Cells.Item(99, 1).Value = 44241
$ws.Cells.Item(99, 2).Value = 0
$ws.Cells.Item(99, 3).Value = 4
$ws.Cells.Item(99, 4).Value = 125.0390747108471
$ws.Cells.Item(100, 1).Value = 44242
$ws.Cells.Item(100, 2).Value = 0
$ws.Cells.Item(100, 3).Value = 5
$ws.Cells.Item(100, 4).Value = 156.2988433885589
$ws.Cells.Item(101, 1).Value = 44243
$ws.Cells.Item(101, 2).Value = 0
$ws.Cells.Item(101, 3).Value = 3
$ws.Cells.Item(101, 4).Value = 93.77930603313536
$ws.Cells.Item(102, 1).Value = 44244
$ws.Cells.Item(102, 2).Value = 0
$ws.Cells.Item(102, 3).Value = 1
$ws.Cells.Item(102, 4).Value = 31.25976867771178
$ws.Cells.Item(103, 1).Value = 44245
$ws.Cells.Item(103, 2).Value = 1
$ws.Cells.Item(103, 3).Value = 5
$ws.Cells.Item(103, 4).Value = 156.2988433885589
$ws.Cells.Item(104, 1).Value = 44246
$ws.Cells.Item(104, 2).Value = 0
$ws.Cells.Item(104, 3).Value = 5
$ws.Cells.Item(104, 4).Value = 156.2988433885589
$ws.Cells.Item(105, 1).Value = 44247
$ws.Cells.Item(105, 2).Value = 0
$ws.Cells.Item(105, 3).Value = 9
$ws.Cells.Item(105, 4).Value = 281.3379180994061
$ws.Cells.Item(106, 1).Value = 44248
$ws.Cells.Item(106, 2).Value = 4
$ws.Cells.Item(106, 3).Value = 9
$ws.Cells.Item(106, 4).Value = 281.3379180994061
$ws.Cells.Item(107, 1).Value = 44249
$ws.Cells.Item(107, 2).Value = 0
$ws.Cells.Item(107, 3).Value = 8
$ws.Cells.Item(107, 4).Value = 250.0781494216943
$ws.Cells.Item(108, 1).Value = 44250
$ws.Cells.Item(108, 2).Value = 4
$ws.Cells.Item(108, 3).Value = 8
$ws.Cells.Item(108, 4).Value = 250.0781494216943
$ws.Cells.Item(109, 1).Value = 44251
$ws.Cells.Item(109, 2).Value = 0
$ws.Cells.Item(109, 3).Value = 9
$ws.Cells.Item(109, 4).Value = 281.3379180994061
$ws.Cells.Item(110, 1).Value = 44252
$ws.Cells.Item(110, 2).Value = 0
$ws.Cells.Item(110, 3).Value = 6
$ws.Cells.Item(110, 4).Value = 187.5586120662707
$ws.Cells.Item(111, 1).Value = 44253
$ws.Cells.Item(111, 2).Value = 0
$ws.Cells.Item(111, 3).Value = 8
$ws.Cells.Item(111, 4).Value = 250.0781494216943
$ws.Cells.Item(112, 1).Value = 44254
$ws.Cells.Item(112, 2).Value = 1
$ws.Cells.Item(112, 3).Value = 10
$ws.Cells.Item(112, 4).Value = 312.5976867771179
$ws.Cells.Item(113, 1).Value = 44255
$ws.Cells.Item(113, 2).Value = 1
$ws.Cells.Item(113, 3).Value = ""
$ws.Cells.Item(113, 4).Value = ""
$ws.Cells.Item(114, 1).Value = 44256
$ws.Cells.Item(114, 2).Value = 2
$ws.Cells.Item(114, 3).Value = ""
$ws.Cells.Item(114, 4).Value = ""
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 6
$ws.Cells.Item(115, 3).Value = ""
$ws.Cells.Item(115, 4).Value = ""

